# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-11-13 (45243) to 2023-11-14 (45244).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) {
    $lastRow = 23
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45243) {
        $cell.Value2 = 45244
    }
}
